$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gross")
$ws.Range("I30").Value = 243
